$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-05-11 Sunday"; new = "2025-05-12 Monday"},
    @{old = "147÷3=49, 0"; new = "606÷9=67, 3"},
    @{old = "714÷2=357, 0"; new = "657÷8=82, 1"},
    @{old = "265÷7=37, 6"; new = "735÷3=245, 0"},
    @{old = "670÷7=95, 5"; new = "639÷7=91, 2"},
    @{old = "598÷9=66, 4"; new = "532÷9=59, 1"},
    @{old = "401÷5=80, 1"; new = "201÷5=40, 1"},
    @{old = "592÷2=296, 0"; new = "916÷2=458, 0"},
    @{old = "165÷5=33, 0"; new = "322÷8=40, 2"},
    @{old = "206÷7=29, 3"; new = "602÷8=75, 2"},
    @{old = "342÷8=42, 6"; new = "209÷5=41, 4"},
    @{old = "734÷4=183, 2"; new = "649÷5=129, 4"},
    @{old = "407÷8=50, 7"; new = "503÷5=100, 3"},
    @{old = "273÷2=136, 1"; new = "671÷9=74, 5"},
    @{old = "503÷9=55, 8"; new = "554÷2=277, 0"},
    @{old = "809÷8=101, 1"; new = "140÷9=15, 5"},
    @{old = "162÷9=18, 0"; new = "371÷2=185, 1"},
    @{old = "868÷4=217, 0"; new = "577÷4=144, 1"},
    @{old = "839÷7=119, 6"; new = "748÷6=124, 4"},
    @{old = "267÷3=89, 0"; new = "699÷6=116, 3"},
    @{old = "963÷5=192, 3"; new = "686÷8=85, 6"},
    @{old = "418÷6=69, 4"; new = "829÷7=118, 3"},
    @{old = "258÷7=36, 6"; new = "922÷9=102, 4"},
    @{old = "178÷5=35, 3"; new = "474÷6=79, 0"},
    @{old = "400÷8=50, 0"; new = "594÷6=99, 0"},
    @{old = "411÷4=102, 3"; new = "266÷5=53, 1"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
